# Transformando em outras classes
$wb = $excel.ActiveWorkbook

# Sheet "Produtos": change F4/G4 from text to numeric values and add G2
$wsProdutos = $wb.Worksheets.Item("Produtos")
$wsProdutos.Range("G2").Value = "nan"
$wsProdutos.Range("F4").Value = 40
$wsProdutos.Range("G4").Value = 20

# Sheet "Estoque": add D2
$wsEstoque = $wb.Worksheets.Item("Estoque")
$wsEstoque.Range("D2").Value = "nan"
